$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Gamelogic" resource-id column (Atlas_ResID) is joined by a second
# config-reference column that ties each row to its SLG building config.
# Row 2 (EFT_INFO) points at the msg_icon config; every other effect row
# shares the common building-setting config ("Ssetting").
$ws.Range("B1").Value = "Atlas_ResID"
$ws.Range("B2").Value = "msg_icon"
$ws.Range("B3:B15").Value = "Ssetting"

[void]$ws.Range("E14").Select()
